$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Heng long Sprocket wheel
$ws.Range("C6").Value = 55
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Heng long Sprocket wheel 1/16 RC Germany Tiger"

# Row 7: Planetary gearbox
$ws.Range("C7").Value = 74.8
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Planetary gearbox for 540 Brushed Motor"

# Update selection to E7
$ws.Range("E7").Select()
